# The document carries two distinct logo pictures that live in the page
# headers/footers:
#   - the Pearson/Edexcel logo (alt text ends in "PearsonLogo.png"), shown
#     twice (first-page footer and default footer), currently named
#     "image2.png" and needs to become "image1.png"
#   - the BTec logo (alt text "BTec_Logo-Orange"), shown once (first-page
#     header), currently named "image1.jpg" and needs to become "image2.jpg"
#
# Rename each inline picture via its InlineShape.Name (the Word object
# model's picture/shape "Name" property, i.e. wp:docPr/@name). The
# AlternativeText (wp:docPr/@descr) is stable/unique per logo, so it is
# used to identify which picture is which before renaming it.

$d = $word.ActiveDocument
$sec = $d.Sections(1)

function Rename-LogoPictures($story) {
    if ($story.Exists) {
        $shapes = $story.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shp = $shapes.Item($j)
            $alt = $shp.AlternativeText
            if ($alt -eq "BTec_Logo-Orange") {
                $shp.Name = "image2.jpg"
            } elseif ($alt -like "*PearsonLogo.png") {
                $shp.Name = "image1.png"
            }
        }
    }
}

for ($hi = 1; $hi -le 3; $hi++) {
    Rename-LogoPictures $sec.Headers($hi)
}
for ($fi = 1; $fi -le 3; $fi++) {
    Rename-LogoPictures $sec.Footers($fi)
}
